$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 45000
$ws.Range("J81").Value = 45000
$ws.Range("L81").Value = 45000
$ws.Range("N81").Value = -46996
$ws.Range("H84").Value = 45000
$ws.Range("J84").Value = 45000
$ws.Range("L84").Value = 135000
$ws.Range("N84").Value = -144984
$ws.Range("H132").Value = 40164110
$ws.Range("I132").Value = 55780344
$ws.Range("J132").Value = 8085.4287
$ws.Range("K132").Value = 167341032
$ws.Range("L132").Value = 24256.2861
$ws.Range("M132").Value = -167338502
$ws.Range("N132").Value = -29316.2861
$ws.Range("H135").Value = 1392.5264
$ws.Range("J135").Value = 3406
$ws.Range("L135").Value = 30654
$ws.Range("N135").Value = -35724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4311.154
$ws.Range("I32").Value = 4290.8105
$ws.Range("K32").Value = 4290.8105
$ws.Range("M32").Value = -4003.8105
$ws.Range("H45").Value = 3478.4614
$ws.Range("I45").Value = 3103
$ws.Range("J45").Value = 3916.5
$ws.Range("K45").Value = 3103
$ws.Range("L45").Value = 3916.5
$ws.Range("M45").Value = -2726
$ws.Range("N45").Value = -4670.5
$ws.Range("H61").Value = 1345.5454
$ws.Range("I61").Value = 1289
$ws.Range("J61").Value = 1600
$ws.Range("K61").Value = 1289
$ws.Range("L61").Value = 1600
$ws.Range("M61").Value = -1077
$ws.Range("N61").Value = -2024
$ws.Range("H74").Value = 297348.34
$ws.Range("I74").Value = 666907.25
$ws.Range("J74").Value = 1701.2
$ws.Range("K74").Value = 666907.25
$ws.Range("L74").Value = 1701.2
$ws.Range("M74").Value = -666033.25
$ws.Range("N74").Value = -3449.2
$ws.Range("H77").Value = 297348.34
$ws.Range("I77").Value = 666907.25
$ws.Range("J77").Value = 1701.2
$ws.Range("K77").Value = 3334536.25
$ws.Range("L77").Value = 8506
$ws.Range("M77").Value = -3330168.25
$ws.Range("N77").Value = -17242
$ws.Range("H136").Value = 1345.5454
$ws.Range("I136").Value = 1289
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 3867
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = -1317
$ws.Range("N136").Value = -9900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 56260
$ws.Range("J59").Value = 56260
$ws.Range("L59").Value = 56260
$ws.Range("N59").Value = -57954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 200917.36
$ws.Range("I31").Value = 376314.84
$ws.Range("J31").Value = 3595.1875
$ws.Range("K31").Value = 376314.84
$ws.Range("L31").Value = 3595.1875
$ws.Range("M31").Value = -376019.84
$ws.Range("N31").Value = -4185.1875
$ws.Range("H34").Value = 200917.36
$ws.Range("I34").Value = 376314.84
$ws.Range("J34").Value = 3595.1875
$ws.Range("K34").Value = 376314.84
$ws.Range("L34").Value = 3595.1875
$ws.Range("M34").Value = -376112.84
$ws.Range("N34").Value = -3999.1875
$ws.Range("H132").Value = 4526.381
$ws.Range("I132").Value = 3482.8
$ws.Range("K132").Value = 10448.4
$ws.Range("M132").Value = -7918.400000000001
$ws.Range("H134").Value = 13823.556
$ws.Range("I134").Value = 13676.5
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 41029.5
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -38494.5
$ws.Range("N134").Value = -50070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3340.4902
$ws.Range("I68").Value = 1495
$ws.Range("J68").Value = 3683.8372
$ws.Range("K68").Value = 4485
$ws.Range("L68").Value = 11051.5116
$ws.Range("M68").Value = -3674
$ws.Range("N68").Value = -12673.5116
$ws.Range("H71").Value = 3340.4902
$ws.Range("I71").Value = 1495
$ws.Range("J71").Value = 3683.8372
$ws.Range("K71").Value = 13455
$ws.Range("L71").Value = 33154.5348
$ws.Range("M71").Value = -9399
$ws.Range("N71").Value = -41266.5348
$ws.Range("H112").Value = 9500
$ws.Range("J112").Value = 9500
$ws.Range("L112").Value = 28500
$ws.Range("N112").Value = -30716
$ws.Range("H113").Value = 461.59183
$ws.Range("J113").Value = 458.89474
$ws.Range("L113").Value = 1376.68422
$ws.Range("N113").Value = -5716.68422
$ws.Range("H131").Value = 758.8125
$ws.Range("I131").Value = 445.2
$ws.Range("J131").Value = 795.27905
$ws.Range("K131").Value = 1335.6
$ws.Range("L131").Value = 2385.83715
$ws.Range("M131").Value = 3704.4
$ws.Range("N131").Value = -12465.83715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2210.75
$ws.Range("I102").Value = 1787.2858
$ws.Range("K102").Value = 1787.2858
$ws.Range("M102").Value = -165.2858000000001
$ws.Range("H132").Value = 4347.343
$ws.Range("I132").Value = 3392.158
$ws.Range("J132").Value = 5481.625
$ws.Range("K132").Value = 10176.474
$ws.Range("L132").Value = 16444.875
$ws.Range("M132").Value = -7646.474
$ws.Range("N132").Value = -21504.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 966.6799999999999
$ws.Range("I68").Value = 924.6301
$ws.Range("J68").Value = 2501.5
$ws.Range("K68").Value = 924.6301
$ws.Range("L68").Value = 2501.5
$ws.Range("M68").Value = -175.6301
$ws.Range("N68").Value = -3999.5
$ws.Range("H71").Value = 966.6799999999999
$ws.Range("I71").Value = 924.6301
$ws.Range("J71").Value = 2501.5
$ws.Range("K71").Value = 4623.1505
$ws.Range("L71").Value = 12507.5
$ws.Range("M71").Value = -879.1504999999997
$ws.Range("N71").Value = -19995.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1495.8572
$ws.Range("I81").Value = 1076.6666
$ws.Range("J81").Value = 2250.4
$ws.Range("K81").Value = 2153.3332
$ws.Range("L81").Value = 4500.8
$ws.Range("M81").Value = -1092.3332
$ws.Range("N81").Value = -6622.8
$ws.Range("H84").Value = 1495.8572
$ws.Range("I84").Value = 1076.6666
$ws.Range("J84").Value = 2250.4
$ws.Range("K84").Value = 10766.666
$ws.Range("L84").Value = 22504
$ws.Range("M84").Value = -5462.666000000001
$ws.Range("N84").Value = -33112
$ws.Range("H126").Value = 562504.75
$ws.Range("I126").Value = 2200.5
$ws.Range("J126").Value = 969998.75
$ws.Range("K126").Value = 6601.5
$ws.Range("L126").Value = 2909996.25
$ws.Range("M126").Value = -4131.5
$ws.Range("N126").Value = -2914936.25
